$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Copy formatting from the last existing data row (639) down through the new rows
$ws.Range("A639:C639").Copy()
$ws.Range("A640:C668").PasteSpecial(-4122)

$ws.Cells.Item(640,1).Value = 'cs'
$ws.Cells.Item(640,2).Value = 'lab.coil.title'
$ws.Cells.Item(640,3).Value = 'Spirálky'

$ws.Cells.Item(641,1).Value = 'cs'
$ws.Cells.Item(641,2).Value = 'lab.coil.subtitle'
$ws.Cells.Item(641,3).Value = 'V této sekci se nachází správa veškerých spirálek, které jste kdy vytvořili.'

$ws.Cells.Item(642,1).Value = 'cs'
$ws.Cells.Item(642,2).Value = 'lab.coil.button.create'
$ws.Cells.Item(642,3).Value = 'Nová spirálka'

$ws.Cells.Item(643,1).Value = 'cs'
$ws.Cells.Item(643,2).Value = 'lab.coil.button.list'
$ws.Cells.Item(643,3).Value = 'Seznam spirálek'

$ws.Cells.Item(644,1).Value = 'cs'
$ws.Cells.Item(644,2).Value = 'lab.coil.list.title'
$ws.Cells.Item(644,3).Value = 'Seznam spirálek'

$ws.Cells.Item(645,1).Value = 'cs'
$ws.Cells.Item(645,2).Value = 'lab.coil.table.wire'
$ws.Cells.Item(645,3).Value = 'Drát'

$ws.Cells.Item(646,1).Value = 'cs'
$ws.Cells.Item(646,2).Value = 'lab.coil.table.wraps'
$ws.Cells.Item(646,3).Value = 'Počet otoček'

$ws.Cells.Item(647,1).Value = 'cs'
$ws.Cells.Item(647,2).Value = 'lab.coil.table.ohm'
$ws.Cells.Item(647,3).Value = 'Odpor'

$ws.Cells.Item(648,1).Value = 'cs'
$ws.Cells.Item(648,2).Value = 'lab.coil.filter.title'
$ws.Cells.Item(648,3).Value = 'Filtrovat spirálky'

$ws.Cells.Item(649,1).Value = 'cs'
$ws.Cells.Item(649,2).Value = 'lab.coil.wireId.label'
$ws.Cells.Item(649,3).Value = 'Drát'

$ws.Cells.Item(650,1).Value = 'cs'
$ws.Cells.Item(650,2).Value = 'lab.coil.preview'
$ws.Cells.Item(650,3).Value = 'Náhled spirálky'

$ws.Cells.Item(651,1).Value = 'cs'
$ws.Cells.Item(651,2).Value = 'lab.coil.preview.preview.title'
$ws.Cells.Item(651,3).Value = 'Detail spirálky'

$ws.Cells.Item(652,1).Value = 'cs'
$ws.Cells.Item(652,2).Value = 'lab.coil.preview.preview.subtitle'
$ws.Cells.Item(652,3).Value = 'Přehled dostupných dat o vybrané spirálce.'

$ws.Cells.Item(653,1).Value = 'cs'
$ws.Cells.Item(653,2).Value = 'lab.coil.preview'
$ws.Cells.Item(653,3).Value = 'Náhled spirálky'

$ws.Cells.Item(654,1).Value = 'cs'
$ws.Cells.Item(654,2).Value = 'lab.coil.button.clone'
$ws.Cells.Item(654,3).Value = 'Klonovat spirálku'

$ws.Cells.Item(655,1).Value = 'cs'
$ws.Cells.Item(655,2).Value = 'lab.coil.button.edit'
$ws.Cells.Item(655,3).Value = 'Upravit spirálku'

$ws.Cells.Item(656,1).Value = 'cs'
$ws.Cells.Item(656,2).Value = 'lab.coil.button.index'
$ws.Cells.Item(656,3).Value = 'Detail spirálky'

$ws.Cells.Item(657,1).Value = 'cs'
$ws.Cells.Item(657,2).Value = 'lab.coil.preview.wire'
$ws.Cells.Item(657,3).Value = 'Drát'

$ws.Cells.Item(658,1).Value = 'cs'
$ws.Cells.Item(658,2).Value = 'lab.coil.preview.wraps'
$ws.Cells.Item(658,3).Value = 'Počet otoček'

$ws.Cells.Item(659,1).Value = 'cs'
$ws.Cells.Item(659,2).Value = 'lab.coil.preview.ohm'
$ws.Cells.Item(659,3).Value = 'Odpor'

$ws.Cells.Item(660,1).Value = 'cs'
$ws.Cells.Item(660,2).Value = 'lab.coil.edit.title'
$ws.Cells.Item(660,3).Value = 'Editace spirálky'

$ws.Cells.Item(661,1).Value = 'cs'
$ws.Cells.Item(661,2).Value = 'lab.coil.edit.subtitle'
$ws.Cells.Item(661,3).Value = 'Upravte vlastnosti vybrané spirálky.'

$ws.Cells.Item(662,1).Value = 'cs'
$ws.Cells.Item(662,2).Value = 'lab.coil.update.submit'
$ws.Cells.Item(662,3).Value = 'Aktualizovat'

$ws.Cells.Item(663,1).Value = 'cs'
$ws.Cells.Item(663,2).Value = 'lab.coil.update.message'
$ws.Cells.Item(663,3).Value = 'Spirálka byla úspěšně aktualizována.'

$ws.Cells.Item(664,1).Value = 'cs'
$ws.Cells.Item(664,2).Value = 'lab.coil.index.title'
$ws.Cells.Item(664,3).Value = 'Detail spirálky'

$ws.Cells.Item(665,1).Value = 'cs'
$ws.Cells.Item(665,2).Value = 'lab.coil.index.preview.subtitle'
$ws.Cells.Item(665,3).Value = 'Přehled dostupných dat o vybrané spirálce.'

$ws.Cells.Item(666,1).Value = 'cs'
$ws.Cells.Item(666,2).Value = 'lab.coil.index.preview.title'
$ws.Cells.Item(666,3).Value = 'Detail spirálky'

$ws.Cells.Item(667,1).Value = 'cs'
$ws.Cells.Item(667,2).Value = 'lab.coil.clone.title'
$ws.Cells.Item(667,3).Value = 'Klon spirálky'

$ws.Cells.Item(668,1).Value = 'cs'
$ws.Cells.Item(668,2).Value = 'lab.coil.link.button'
$ws.Cells.Item(668,3).Value = 'Detail spirálky'

# Restore view state: scroll position and active selection
$ws.Range("B658").Select()
